$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 7 de Abril de 2020 a las 13:22"

# Row 5
$ws.Range("C5").Value = 10738
$ws.Range("D5").Value = 14677

# Row 8
$ws.Range("B8").Value = 4298
$ws.Range("C8").Value = 1353
$ws.Range("D8").Value = 8547
$ws.Range("E8").Value = 370

# Row 12 (Albacete)
$ws.Range("A12").Value = "Albacete"
$ws.Range("B12").Value = 2832
$ws.Range("C12").Value = 1353
$ws.Range("D12").Value = 8547
$ws.Range("E12").Value = 252

# Row 13 (Castilla-La Mancha)
$ws.Range("A13").Value = "Castilla-La Mancha"
$ws.Range("B13").Value = 2780
$ws.Range("C13").Value = 71
$ws.Range("D13").Value = 2446
$ws.Range("E13").Value = 263

# Row 17
$ws.Range("B17").Value = 2434
$ws.Range("C17").Value = 1353
$ws.Range("D17").Value = 8547
$ws.Range("E17").Value = 331

# Row 32
$ws.Range("C32").Value = 192
$ws.Range("D32").Value = 1449

# Row 37 (Guadalajara)
$ws.Range("A37").Value = "Guadalajara"
$ws.Range("B37").Value = 897
$ws.Range("C37").Value = 1353
$ws.Range("D37").Value = 8547
$ws.Range("E37").Value = 128

# Row 38 (Castello/Castellon)
$ws.Range("A38").Value = "Castello/Castellon"
$ws.Range("B38").Value = 876
$ws.Range("C38").Value = 107
$ws.Range("D38").Value = 691
$ws.Range("E38").Value = 78

# Row 44 (Cuenca)
$ws.Range("A44").Value = "Cuenca"
$ws.Range("B44").Value = 616
$ws.Range("C44").Value = 1353
$ws.Range("D44").Value = 8547
$ws.Range("E44").Value = 96

# Row 45 (Lugo)
$ws.Range("A45").Value = "Lugo"
$ws.Range("B45").Value = 586
$ws.Range("C45").Value = 333
$ws.Range("D45").Value = 520
$ws.Range("E45").Value = 11

# Row 48
$ws.Range("C48").Value = 192
$ws.Range("D48").Value = 1449

# Row 55
$ws.Range("B55").Value = 84
$ws.Range("C55").Value = 7
$ws.Range("D55").Value = 73

# Row 56
$ws.Range("C56").Value = 192
$ws.Range("D56").Value = 1449

# Row 57
$ws.Range("C57").Value = 192
$ws.Range("D57").Value = 1449

# Row 59
$ws.Range("C59").Value = 192
$ws.Range("D59").Value = 1449

# Row 62
$ws.Range("C62").Value = 192
$ws.Range("D62").Value = 1449

# Row 64
$ws.Range("C64").Value = 192
$ws.Range("D64").Value = 1449
